# Auto-generated edit script for 7_gunluk workbook
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update header text: "Random Forest" -> "Decision Tree"
$ws.Range("C1").Value = "Decision Tree"

# Update changed numeric values in columns C (predicted) and D (actual)
$ws.Range("C12").Value = 2.1
$ws.Range("D12").Value = 2
$ws.Range("C13").Value = 2.169230769230769
$ws.Range("D13").Value = 2
$ws.Range("C14").Value = 3.2
$ws.Range("D14").Value = 3
$ws.Range("C15").Value = 4.75
$ws.Range("D15").Value = 4
$ws.Range("C16").Value = 6.36
$ws.Range("D16").Value = 7
$ws.Range("C17").Value = 6.279999999999999
$ws.Range("D17").Value = 7
$ws.Range("C18").Value = 5.892307692307692
$ws.Range("D18").Value = 6
$ws.Range("C19").Value = 5.957142857142857
$ws.Range("D19").Value = 6
$ws.Range("C20").Value = 5.316666666666666
$ws.Range("D20").Value = 5
$ws.Range("C21").Value = 5.177777777777778
$ws.Range("D21").Value = 5
$ws.Range("C22").Value = 4.933333333333334
$ws.Range("D22").Value = 4
$ws.Range("C23").Value = 3.175
$ws.Range("D23").Value = 3
$ws.Range("C24").Value = 0.2363636363636363
$ws.Range("D24").Value = 0
$ws.Range("C25").Value = 0
$ws.Range("C36").Value = 1
$ws.Range("D36").Value = 1
$ws.Range("C37").Value = 1
$ws.Range("D37").Value = 1
$ws.Range("C38").Value = 3.38
$ws.Range("D38").Value = 3
$ws.Range("C39").Value = 3.523076923076923
$ws.Range("D39").Value = 3
$ws.Range("C40").Value = 5.55
$ws.Range("D40").Value = 6
$ws.Range("C41").Value = 5.781818181818182
$ws.Range("D41").Value = 6
$ws.Range("C42").Value = 5.8
$ws.Range("D42").Value = 6
$ws.Range("C43").Value = 5.533333333333333
$ws.Range("D43").Value = 5
$ws.Range("C44").Value = 5.422222222222222
$ws.Range("D44").Value = 5
$ws.Range("C45").Value = 5.279999999999999
$ws.Range("D45").Value = 5
$ws.Range("C46").Value = 5.173333333333333
$ws.Range("D46").Value = 5
$ws.Range("C47").Value = 3.381818181818182
$ws.Range("D47").Value = 3
$ws.Range("C48").Value = 0.05454545454545454
$ws.Range("D48").Value = 0
$ws.Range("C49").Value = 0
$ws.Range("C60").Value = 1.757142857142857
$ws.Range("D60").Value = 2
$ws.Range("C61").Value = 2.2
$ws.Range("D61").Value = 2
$ws.Range("C62").Value = 3.311111111111111
$ws.Range("D62").Value = 3
$ws.Range("C63").Value = 4.357142857142857
$ws.Range("D63").Value = 4
$ws.Range("C64").Value = 6.266666666666667
$ws.Range("D64").Value = 7
$ws.Range("C65").Value = 6.185714285714285
$ws.Range("D65").Value = 7
$ws.Range("C66").Value = 5.857142857142857
$ws.Range("D66").Value = 6
$ws.Range("C67").Value = 5.6625
$ws.Range("D67").Value = 6
$ws.Range("C68").Value = 5.476923076923077
$ws.Range("D68").Value = 5
$ws.Range("C69").Value = 5.166666666666667
$ws.Range("D69").Value = 5
$ws.Range("C70").Value = 4.716666666666667
$ws.Range("D70").Value = 4
$ws.Range("C71").Value = 2.98
$ws.Range("D71").Value = 3
$ws.Range("C72").Value = 0.2285714285714286
$ws.Range("D72").Value = 0
$ws.Range("C73").Value = 0
$ws.Range("D73").Value = 0
$ws.Range("C81").Value = 0.08888888888888889
$ws.Range("C84").Value = 2.127272727272727
$ws.Range("D84").Value = 2
$ws.Range("C85").Value = 2.383333333333333
$ws.Range("D85").Value = 2
$ws.Range("C86").Value = 4.08
$ws.Range("D86").Value = 4
$ws.Range("C87").Value = 4.566666666666666
$ws.Range("D87").Value = 4
$ws.Range("C88").Value = 6.15
$ws.Range("D88").Value = 7
$ws.Range("C89").Value = 6.43076923076923
$ws.Range("D89").Value = 7
$ws.Range("C90").Value = 6.02
$ws.Range("D90").Value = 6
$ws.Range("C91").Value = 5.969230769230768
$ws.Range("D91").Value = 6
$ws.Range("C92").Value = 6.24
$ws.Range("D92").Value = 6
$ws.Range("C93").Value = 5.4
$ws.Range("D93").Value = 5
$ws.Range("C94").Value = 4.733333333333333
$ws.Range("D94").Value = 4
$ws.Range("C95").Value = 3.56
$ws.Range("D95").Value = 3
$ws.Range("C96").Value = 0.2
$ws.Range("D96").Value = 0
$ws.Range("C97").Value = 0
$ws.Range("C108").Value = 2.581818181818182
$ws.Range("D108").Value = 3
$ws.Range("C109").Value = 3.16
$ws.Range("D109").Value = 3
$ws.Range("C110").Value = 4.333333333333333
$ws.Range("D110").Value = 4
$ws.Range("C111").Value = 5.272727272727272
$ws.Range("D111").Value = 5
$ws.Range("C112").Value = 7.018181818181819
$ws.Range("D112").Value = 8
$ws.Range("C113").Value = 7.109090909090909
$ws.Range("D113").Value = 8
$ws.Range("C114").Value = 7.38
$ws.Range("D114").Value = 7
$ws.Range("C115").Value = 7.018181818181819
$ws.Range("D115").Value = 7
$ws.Range("C116").Value = 6.506666666666666
$ws.Range("D116").Value = 7
$ws.Range("C117").Value = 5.057142857142857
$ws.Range("D117").Value = 5
$ws.Range("C118").Value = 4.246153846153846
$ws.Range("D118").Value = 4
$ws.Range("C119").Value = 2.422222222222222
$ws.Range("D119").Value = 3
$ws.Range("C120").Value = 0.38
$ws.Range("D120").Value = 0
$ws.Range("C121").Value = 0
$ws.Range("C129").Value = 0.02
$ws.Range("C132").Value = 1.085714285714286
$ws.Range("D132").Value = 1
$ws.Range("C133").Value = 0.9714285714285714
$ws.Range("D133").Value = 1
$ws.Range("C134").Value = 3.377777777777778
$ws.Range("D134").Value = 3
$ws.Range("C135").Value = 3.78
$ws.Range("D135").Value = 3
$ws.Range("C136").Value = 5.399999999999999
$ws.Range("D136").Value = 6
$ws.Range("C137").Value = 5.48
$ws.Range("D137").Value = 6
$ws.Range("C138").Value = 5.292307692307692
$ws.Range("D138").Value = 5
$ws.Range("C139").Value = 5.063157894736842
$ws.Range("D139").Value = 5
$ws.Range("C140").Value = 5
$ws.Range("D140").Value = 5
$ws.Range("C141").Value = 5.090909090909091
$ws.Range("D141").Value = 5
$ws.Range("C142").Value = 4.2
$ws.Range("D142").Value = 4
$ws.Range("C143").Value = 2.76
$ws.Range("D143").Value = 3
$ws.Range("C144").Value = 0.0380952380952381
$ws.Range("D144").Value = 0
$ws.Range("C145").Value = 0
$ws.Range("C153").Value = 0.04
$ws.Range("C156").Value = 2.8125
$ws.Range("D156").Value = 3
$ws.Range("C157").Value = 3.4
$ws.Range("D157").Value = 3
$ws.Range("C158").Value = 4.279999999999999
$ws.Range("D158").Value = 4
$ws.Range("C159").Value = 6.2
$ws.Range("D159").Value = 5
$ws.Range("C160").Value = 6.85
$ws.Range("D160").Value = 8
$ws.Range("C161").Value = 7.145454545454545
$ws.Range("D161").Value = 8
$ws.Range("C162").Value = 7.199999999999999
$ws.Range("D162").Value = 8
$ws.Range("C163").Value = 7.018181818181819
$ws.Range("D163").Value = 7
$ws.Range("C164").Value = 7.12
$ws.Range("D164").Value = 7
$ws.Range("C165").Value = 5.892307692307692
$ws.Range("D165").Value = 5
$ws.Range("C166").Value = 4.6
$ws.Range("D166").Value = 4
$ws.Range("C167").Value = 3.45
$ws.Range("D167").Value = 3
$ws.Range("C168").Value = 0.15
$ws.Range("D168").Value = 0
$ws.Range("C169").Value = 0.08
